$wb = $excel.ActiveWorkbook

# --- Data for the two new sheets -------------------------------------------------
$dataPersistence = @(
@("Event", "Lead Time 1", "Lead Time 2", "Lead Time 3", "Lead Time 4", "Lead Time 5", "Lead Time 6"),
@("20230520_2235", 0.1093417941395427, 0.14222355786868809, 0.15801002784225149, 0.16965018884312891, 0.17830929515869909, 0.18457944185765851),
@("20190320_0005", 0.089383744817873784, 0.1213976531881588, 0.1384881897847427, 0.14976766349115139, 0.1579584453755801, 0.16409101931762449),
@("20191222_0900", 0.14349619854771259, 0.1486672312352991, 0.1468072376199161, 0.1436444739009787, 0.1379804444130428, 0.13623532815140549),
@("20180505_1745", 0.62705033708427094, 0.77769100920172141, 0.85315357485029975, 0.89594382561202379, 0.91973496241805364, 0.92998444516056111),
@("20230513_1455", 0.28856055506356321, 0.28261419398977128, 0.28631599365266741, 0.28183186357125001, 0.28381551123303, 0.29355704340687222),
@("20200911_1315", 1.1391007590020079, 1.544924529805995, 1.7517193610172239, 1.8681238748736471, 1.936648119795648, 1.965836376779575),
@("20191111_0710", 0.15347895665286729, 0.17396842408641791, 0.1795493887662763, 0.1832357416743266, 0.18036340927272049, 0.1811730361559942),
@("20230302_0245", 0.12714267696041601, 0.15467071207716759, 0.14824028500771269, 0.14311314727084751, 0.1464842130220288, 0.15131245666838691),
@("20190412_1220", 0.2679306195213278, 0.34117741846858352, 0.37595133209369669, 0.40741503767713771, 0.43335559780736288, 0.44982570696938018),
@("20200120_1440", 0.14794309867913441, 0.1585449592400644, 0.1590537357547557, 0.15911507679094211, 0.1605507893470518, 0.15929201541795779),
@("20230129_2215", 0.039289593368064743, 0.05148791094741352, 0.059545184691728217, 0.064132073702365233, 0.06644363645283724, 0.067493813184186041),
@("20181014_0515", 2.6818337324741641, 3.312149919344368, 3.5536470993431819, 3.6934733795820311, 3.7631099364032421, 3.782436419408485)
)

$dataExtrapolation = @(
@("Event", "Lead Time 1", "Lead Time 2", "Lead Time 3", "Lead Time 4", "Lead Time 5", "Lead Time 6"),
@("20230520_2235", 0.072248601410728738, 0.1077421647809837, 0.1326907574350096, 0.1533659477448128, 0.17097806205198199, 0.18596631786870849),
@("20190320_0005", 0.054729302274071918, 0.082272740464693447, 0.1023754027054072, 0.1187219453815359, 0.13372162053131781, 0.14723232477136541),
@("20191222_0900", 0.15520230839188179, 0.19775071311360051, 0.14162666436798979, 0.18079457013762251, 0.18228570471910949, 0.20095360582903449),
@("20180505_1745", 0.54178500361030324, 0.77450528464405555, 0.92326472194882037, 1.030046949780848, 1.104603656689586, 1.1517878904351451),
@("20230513_1455", 0.28117932300678727, 0.3038213660094049, 0.29155030490730499, 0.26633476133093831, 0.25431757930790738, 0.1797641632760644),
@("20200911_1315", 1.0219523177301051, 1.495120945455747, 1.7764108756864561, 1.9644666916128981, 2.1330018762000811, 2.2647382398480582),
@("20191111_0710", 0.1051869864806684, 0.1527935918909292, 0.18481116523437149, 0.20742773676623669, 0.22256232256211411, 0.22992672790789589),
@("20230302_0245", 0.098245351551857332, 0.15201673963283141, 0.1872774002338165, 0.20455651076824949, 0.22213365277678501, 0.22976503803190729),
@("20190412_1220", 0.210682903068767, 0.31256785810953752, 0.37461898438198682, 0.42141319356573531, 0.45676677768148233, 0.48166405857396638),
@("20200120_1440", 0.13121265882287039, 0.15206882788484, 0.1761266715043622, 0.1877291704750505, 0.1897271941203367, 0.17180396560569261),
@("20230129_2215", 0.034066695701907831, 0.052534501406709858, 0.066752686723907842, 0.079163099528389416, 0.089063589426906037, 0.096804577793057733),
@("20181014_0515", 2.31901243164, 3.2121306330270718, 3.6860239096763441, 4.0335337658123587, 4.3297292995232404, 4.5626489211574048)
)

# --- Create "VET_persistence" by copying the last existing sheet (keeps header style) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Copy($null, $lastSheet)
$vetPersistence = $wb.Worksheets.Item($wb.Worksheets.Count)
$vetPersistence.Name = "VET_persistence"

for ($r = 0; $r -lt $dataPersistence.Length; $r++) {
    $row = $dataPersistence[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $vetPersistence.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# --- Create "VET_extrapolation" by copying VET_persistence (same header style) ---
$vetPersistence.Copy($null, $vetPersistence)
$vetExtrapolation = $wb.Worksheets.Item($wb.Worksheets.Count)
$vetExtrapolation.Name = "VET_extrapolation"

for ($r = 0; $r -lt $dataExtrapolation.Length; $r++) {
    $row = $dataExtrapolation[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $vetExtrapolation.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# --- Make the newly-added last sheet ("VET_extrapolation") the active tab ---
$vetExtrapolation.Activate()

Write-Host "Sheets now:"
foreach ($s in $wb.Worksheets) { Write-Host $s.Name }
